$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update two recalculated amounts (early repayment figures)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 672.06
$wsSummary.Range("E3").Value = 672.06

# ---------------------------------------------------------------------------
# Sheet "Repayment Schedule": update recalculated interest/total figures for
# row 6, and shift the trailing blank cell in row 2 from column O to column P
# (so it lines up under the "Outstanding" header, like every other row).
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("H6").Value = 73.97
$wsRepay.Range("K6").Value = 907.3
$wsRepay.Range("P6").Value = 907.3

$wsRepay.Range("N2").Copy()
$wsRepay.Range("P2").PasteSpecial(-4122)   # xlPasteFormats
$wsRepay.Range("O2").Clear()

# ---------------------------------------------------------------------------
# Update the saved selections (cursor position) on each sheet. Do this last,
# finishing on "Repayment Schedule" so it remains the active tab, since
# selecting a range on a sheet activates that sheet/tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("NewLoanInput").Range("B3:B15").Select()
$wb.Worksheets.Item("Summary").Range("D17").Select()
$wsRepay.Range("P1:P14").Select()
